$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("M2").Value = 3.442313
$ws.Range("N2").Value = 10.326939
$ws.Range("O2").Value = 0.2120737065114005
$ws.Range("P2").Value = 0.2368526181325179
$ws.Range("Q2").Value = 3.965807339225667
$ws.Range("R2").Value = 35.692266053031
$ws.Range("S2").Value = 0.2120737065114005
$ws.Range("T2").Value = 0.2368526181325179

$ws.Range("O3").Value = 0.158453673516874
$ws.Range("P3").Value = 0.1769675649214407
$ws.Range("S3").Value = 0.158453673516874
$ws.Range("T3").Value = 0.1769675649214407

$ws.Range("M4").Value = 2.535264
$ws.Range("N4").Value = 7.605791999999999
$ws.Range("O4").Value = 0.1561923141402073
$ws.Range("P4").Value = 0.174441985971967
$ws.Range("Q4").Value = 2.920817653152
$ws.Range("R4").Value = 26.287358878368
$ws.Range("S4").Value = 0.1561923141402073
$ws.Range("T4").Value = 0.174441985971967

$ws.Range("M5").Value = 5.0943505
$ws.Range("N5").Value = 10.188701
$ws.Range("O5").Value = 0.3138522826957358
$ws.Range("P5").Value = 0.2336820724146239
$ws.Range("Q5").Value = 5.869080644754833
$ws.Range("R5").Value = 35.214483868529
$ws.Range("S5").Value = 0.3138522826957358
$ws.Range("T5").Value = 0.2336820724146239

$ws.Range("M6").Value = 2.587785
$ws.Range("N6").Value = 7.763355
$ws.Range("O6").Value = 0.1594280231357824
$ws.Range("P6").Value = 0.1780557585594505
$ws.Range("Q6").Value = 2.981325854255
$ws.Range("R6").Value = 26.831932688295
$ws.Range("S6").Value = 0.1594280231357824
$ws.Range("T6").Value = 0.1780557585594505
